$d = $word.ActiveDocument

# --- 1. "Week 8" heading: append date range text (split across a few
#        InsertAfter calls, mirroring how this was typed/edited in Word) ---
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Heading 2" -and $p.Range.Text.TrimEnd([char]13) -eq "Week 8") {
        $r = $d.Range($p.Range.End - 1, $p.Range.End - 1)
        $r.InsertAfter(": Nov 2 – Nov")
        $r2 = $d.Range($r.End, $r.End)
        $r2.InsertAfter(" 8")
        $r3 = $d.Range($r2.End, $r2.End)
        $r3.InsertAfter(" ")
        break
    }
}

# --- 2. Tidy up double spaces in the last Week 8 bullet ---
$d.Content.Find.Execute("append  if", $true, $false, $false, $false, $false, `
    $true, 1, $false, "append if", 2) | Out-Null
$d.Content.Find.Execute("in  multiple", $true, $false, $false, $false, $false, `
    $true, 1, $false, "in multiple", 2) | Out-Null

# --- 3. Append the new Week 9 section: blank spacer line, heading, and
#        two new bulleted list items describing the database class work ---
$lastP = $d.Paragraphs.Last
$insertionPoint = $d.Range($lastP.Range.End, $lastP.Range.End)

$flatOpc = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData>' + `
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr></w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Week 9: Nov 9 – Nov 15</w:t></w:r></w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t>Added text file containing words with difficulty level and hint.</w:t></w:r></w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t xml:space="preserve">Implemented a database class to load the text file into a table and perform CRUD operations on the table.</w:t></w:r></w:p>' + `
  '</w:body></w:document>' + `
  '</pkg:xmlData></pkg:part></pkg:package>'

$insertionPoint.InsertXML($flatOpc) | Out-Null

# --- 4. Turn the two new bullets into a numbered list (this mints a
#        fresh list definition / numId, same as clicking "Numbering" on
#        a brand-new list in Word) ---
$bullet1 = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$bullet2 = $d.Paragraphs.Item($d.Paragraphs.Count)

$template = $d.ListTemplates.Item(1)
$bullet1.Range.ListFormat.ApplyListTemplateWithLevel($template)
$bullet2.Range.ListFormat.ApplyListTemplateWithLevel($template, $true)

Write-Output "Edit complete. Paragraphs: $($d.Paragraphs.Count)"
